$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(149, 1).Value = "6_18"
$ws.Cells.Item(149, 2).Value = 56
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 0
$ws.Cells.Item(149, 5).Value = 0
$ws.Cells.Item(149, 6).Value = 51
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 1
$ws.Cells.Item(149, 9).Value = 0

$ws.Cells.Item(150, 1).Value = "7_00"
$ws.Cells.Item(150, 2).Value = 54
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 0
$ws.Cells.Item(150, 5).Value = 0
$ws.Cells.Item(150, 6).Value = 39
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 1
$ws.Cells.Item(150, 9).Value = 0

$ws.Cells.Item(151, 1).Value = "7_06"
$ws.Cells.Item(151, 2).Value = 8
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 0
$ws.Cells.Item(151, 5).Value = 0
$ws.Cells.Item(151, 6).Value = 40
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 1
$ws.Cells.Item(151, 9).Value = 0

$ws.Cells.Item(152, 1).Value = "7_12"
$ws.Cells.Item(152, 2).Value = 57
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 0
$ws.Cells.Item(152, 5).Value = 0
$ws.Cells.Item(152, 6).Value = 58
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 1
$ws.Cells.Item(152, 9).Value = 0

$ws.Cells.Item(153, 1).Value = "7_18"
$ws.Cells.Item(153, 2).Value = 55
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 0
$ws.Cells.Item(153, 5).Value = 0
$ws.Cells.Item(153, 6).Value = 52
$ws.Cells.Item(153, 7).Value = 0
$ws.Cells.Item(153, 8).Value = 1
$ws.Cells.Item(153, 9).Value = 0

$ws.Cells.Item(154, 1).Value = "8_00"
$ws.Cells.Item(154, 2).Value = 49
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(154, 4).Value = 0
$ws.Cells.Item(154, 5).Value = 0
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 1
$ws.Cells.Item(154, 9).Value = 0

$ws.Cells.Item(155, 1).Value = "8_06"
$ws.Cells.Item(155, 2).Value = 41
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 0
$ws.Cells.Item(155, 5).Value = 0
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 1
$ws.Cells.Item(155, 9).Value = 0

$ws.Cells.Item(156, 1).Value = "8_12"
$ws.Cells.Item(156, 2).Value = 0
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 0
$ws.Cells.Item(156, 5).Value = 0
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 1
$ws.Cells.Item(156, 9).Value = 0

$ws.Cells.Item(157, 1).Value = "8_18"
$ws.Cells.Item(157, 2).Value = 0
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 0
$ws.Cells.Item(157, 5).Value = 0
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 1
$ws.Cells.Item(157, 9).Value = 0

$ws.Cells.Item(158, 1).Value = "9_00"
$ws.Cells.Item(158, 2).Value = 43
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = 0
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 1
$ws.Cells.Item(158, 9).Value = 0

$ws.Cells.Item(159, 1).Value = "9_06"
$ws.Cells.Item(159, 2).Value = 35
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 5).Value = 0
$ws.Cells.Item(159, 6).Value = 25
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 1
$ws.Cells.Item(159, 9).Value = 0

$ws.Cells.Item(160, 1).Value = "9_12"
$ws.Cells.Item(160, 2).Value = 50
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 0
$ws.Cells.Item(160, 5).Value = 0
$ws.Cells.Item(160, 6).Value = 55
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 1
$ws.Cells.Item(160, 9).Value = 0

$ws.Cells.Item(161, 1).Value = "9_18"
$ws.Cells.Item(161, 2).Value = 45
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 0
$ws.Cells.Item(161, 5).Value = 0
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 1
$ws.Cells.Item(161, 9).Value = 0

$ws.Cells.Item(162, 1).Value = "10_00"
$ws.Cells.Item(162, 2).Value = 53
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 0
$ws.Cells.Item(162, 5).Value = 0
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 1
$ws.Cells.Item(162, 9).Value = 0

$ws.Cells.Item(163, 1).Value = "10_06"
$ws.Cells.Item(163, 2).Value = 38
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = 0
$ws.Cells.Item(163, 5).Value = 0
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 1
$ws.Cells.Item(163, 9).Value = 0

$ws.Cells.Item(164, 1).Value = "10_12"
$ws.Cells.Item(164, 2).Value = 61
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 5).Value = 0
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 1
$ws.Cells.Item(164, 9).Value = 0

$ws.Cells.Item(165, 1).Value = "10_18"
$ws.Cells.Item(165, 2).Value = 47
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 5).Value = 0
$ws.Cells.Item(165, 6).Value = 0
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 1
$ws.Cells.Item(165, 9).Value = 0

$ws.Cells.Item(166, 1).Value = "11_00"
$ws.Cells.Item(166, 2).Value = 62
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 0
$ws.Cells.Item(166, 5).Value = 0
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 1
$ws.Cells.Item(166, 9).Value = 0

$ws.Cells.Item(167, 1).Value = "11_06"
$ws.Cells.Item(167, 2).Value = 53
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 0
$ws.Cells.Item(167, 5).Value = 0
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 1
$ws.Cells.Item(167, 9).Value = 0

$ws.Cells.Item(168, 1).Value = "11_12"
$ws.Cells.Item(168, 2).Value = 0
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 0
$ws.Cells.Item(168, 5).Value = 0
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 1
$ws.Cells.Item(168, 8).Value = 98
$ws.Cells.Item(168, 9).Value = 0

$ws.Cells.Item(169, 1).Value = "11_18"
$ws.Cells.Item(169, 2).Value = 0
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(169, 4).Value = 0
$ws.Cells.Item(169, 5).Value = 0
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 1
$ws.Cells.Item(169, 8).Value = 97
$ws.Cells.Item(169, 9).Value = 0

$ws.Cells.Item(170, 1).Value = "12_00"
$ws.Cells.Item(170, 2).Value = 0
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 0
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 1
$ws.Cells.Item(170, 8).Value = 96
$ws.Cells.Item(170, 9).Value = 0

$ws.Cells.Item(171, 1).Value = "12_06"
$ws.Cells.Item(171, 2).Value = 0
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 5).Value = 0
$ws.Cells.Item(171, 6).Value = 1
$ws.Cells.Item(171, 7).Value = 74
$ws.Cells.Item(171, 8).Value = 97
$ws.Cells.Item(171, 9).Value = 0

$ws.Cells.Item(172, 1).Value = "12_12"
$ws.Cells.Item(172, 2).Value = 0
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = 0
$ws.Cells.Item(172, 6).Value = 1
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 98
$ws.Cells.Item(172, 9).Value = 0

$ws.Cells.Item(173, 1).Value = "12_18"
$ws.Cells.Item(173, 2).Value = 0
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 0
$ws.Cells.Item(173, 6).Value = 1
$ws.Cells.Item(173, 7).Value = 78
$ws.Cells.Item(173, 8).Value = 97
$ws.Cells.Item(173, 9).Value = 0

$ws.Cells.Item(174, 1).Value = "13_00"
$ws.Cells.Item(174, 2).Value = 0
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 0
$ws.Cells.Item(174, 6).Value = 1
$ws.Cells.Item(174, 7).Value = 77
$ws.Cells.Item(174, 8).Value = 98
$ws.Cells.Item(174, 9).Value = 0

$ws.Cells.Item(175, 1).Value = "13_06"
$ws.Cells.Item(175, 2).Value = 0
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 5).Value = 0
$ws.Cells.Item(175, 6).Value = 1
$ws.Cells.Item(175, 7).Value = 80
$ws.Cells.Item(175, 8).Value = 96
$ws.Cells.Item(175, 9).Value = 0

$ws.Cells.Item(176, 1).Value = "13_12"
$ws.Cells.Item(176, 2).Value = 0
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 5).Value = 0
$ws.Cells.Item(176, 6).Value = 1
$ws.Cells.Item(176, 7).Value = 80
$ws.Cells.Item(176, 8).Value = 97
$ws.Cells.Item(176, 9).Value = 0

$ws.Cells.Item(177, 1).Value = "13_18"
$ws.Cells.Item(177, 2).Value = 0
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 0
$ws.Cells.Item(177, 6).Value = 1
$ws.Cells.Item(177, 7).Value = 79
$ws.Cells.Item(177, 8).Value = 97
$ws.Cells.Item(177, 9).Value = 0

$ws.Cells.Item(178, 1).Value = "14_00"
$ws.Cells.Item(178, 2).Value = 0
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 0
$ws.Cells.Item(178, 5).Value = 0
$ws.Cells.Item(178, 6).Value = 1
$ws.Cells.Item(178, 7).Value = 49
$ws.Cells.Item(178, 8).Value = 97
$ws.Cells.Item(178, 9).Value = 0

$ws.Cells.Item(179, 1).Value = "14_12"
$ws.Cells.Item(179, 2).Value = 0
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 5).Value = 0
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0
$ws.Cells.Item(179, 9).Value = 0

$ws.Cells.Item(180, 1).Value = "14_18"
$ws.Cells.Item(180, 2).Value = 0
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 0
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0
$ws.Cells.Item(180, 9).Value = 0

$ws.Cells.Item(181, 1).Value = "15_00"
$ws.Cells.Item(181, 2).Value = 0
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 0
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0
$ws.Cells.Item(181, 9).Value = 0

$ws.Cells.Item(182, 1).Value = "15_06"
$ws.Cells.Item(182, 2).Value = 0
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 0
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0
$ws.Cells.Item(182, 9).Value = 0

$ws.Cells.Item(183, 1).Value = "15_12"
$ws.Cells.Item(183, 2).Value = 35
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 0
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0
$ws.Cells.Item(183, 9).Value = 0

$ws.Cells.Item(184, 1).Value = "15_18"
$ws.Cells.Item(184, 2).Value = 52
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 0
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0
$ws.Cells.Item(184, 9).Value = 0

$ws.Cells.Item(185, 1).Value = "16_00"
$ws.Cells.Item(185, 2).Value = 48
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 0
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0
$ws.Cells.Item(185, 9).Value = 0

$ws.Cells.Item(186, 1).Value = "16_06"
$ws.Cells.Item(186, 2).Value = 22
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 0
$ws.Cells.Item(186, 6).Value = 29
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0
$ws.Cells.Item(186, 9).Value = 0

$ws.Cells.Item(187, 1).Value = "16_12"
$ws.Cells.Item(187, 2).Value = 60
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 0
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0
$ws.Cells.Item(187, 9).Value = 0

$ws.Cells.Item(188, 1).Value = "16_18"
$ws.Cells.Item(188, 2).Value = 53
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 5).Value = 0
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0
$ws.Cells.Item(188, 9).Value = 0

$ws.Cells.Item(189, 1).Value = "17_00"
$ws.Cells.Item(189, 2).Value = 61
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 0
$ws.Cells.Item(189, 5).Value = 0
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 0
$ws.Cells.Item(189, 9).Value = 0

$ws.Cells.Item(190, 1).Value = "17_06"
$ws.Cells.Item(190, 2).Value = 40
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 5).Value = 0
$ws.Cells.Item(190, 6).Value = 37
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0
$ws.Cells.Item(190, 9).Value = 0

$ws.Cells.Item(191, 1).Value = "17_12"
$ws.Cells.Item(191, 2).Value = 35
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 0
$ws.Cells.Item(191, 6).Value = 35
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 0
$ws.Cells.Item(191, 9).Value = 0

$ws.Cells.Item(192, 1).Value = "17_18"
$ws.Cells.Item(192, 2).Value = 35
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 5).Value = 0
$ws.Cells.Item(192, 6).Value = 33
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0
$ws.Cells.Item(192, 9).Value = 0

$ws.Cells.Item(193, 1).Value = "18_00"
$ws.Cells.Item(193, 2).Value = 19
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 0
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 27
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0
$ws.Cells.Item(193, 9).Value = 0

$ws.Cells.Item(194, 1).Value = "18_06"
$ws.Cells.Item(194, 2).Value = 46
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 0
$ws.Cells.Item(194, 5).Value = 0
$ws.Cells.Item(194, 6).Value = 33
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0
$ws.Cells.Item(194, 9).Value = 0

$ws.Cells.Item(195, 1).Value = "18_12"
$ws.Cells.Item(195, 2).Value = 35
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 5).Value = 0
$ws.Cells.Item(195, 6).Value = 33
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0
$ws.Cells.Item(195, 9).Value = 0

$ws.Cells.Item(196, 1).Value = "18_18"
$ws.Cells.Item(196, 2).Value = 54
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 0
$ws.Cells.Item(196, 5).Value = 0
$ws.Cells.Item(196, 6).Value = 45
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 0
$ws.Cells.Item(196, 9).Value = 0

$ws.Cells.Item(197, 1).Value = "19_00"
$ws.Cells.Item(197, 2).Value = 35
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 0
$ws.Cells.Item(197, 5).Value = 0
$ws.Cells.Item(197, 6).Value = 33
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 0
$ws.Cells.Item(197, 9).Value = 0

$ws.Cells.Item(198, 1).Value = "19_06"
$ws.Cells.Item(198, 2).Value = 43
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 0
$ws.Cells.Item(198, 5).Value = 0
$ws.Cells.Item(198, 6).Value = 45
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 0
$ws.Cells.Item(198, 9).Value = 0

$ws.Cells.Item(199, 1).Value = "19_12"
$ws.Cells.Item(199, 2).Value = 56
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 0
$ws.Cells.Item(199, 5).Value = 0
$ws.Cells.Item(199, 6).Value = 30
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0
$ws.Cells.Item(199, 9).Value = 0

$ws.Cells.Item(200, 1).Value = "19_18"
$ws.Cells.Item(200, 2).Value = 34
$ws.Cells.Item(200, 3).Value = 0
$ws.Cells.Item(200, 4).Value = 0
$ws.Cells.Item(200, 5).Value = 0
$ws.Cells.Item(200, 6).Value = 30
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = 0
$ws.Cells.Item(200, 9).Value = 0

$ws.Cells.Item(201, 1).Value = "20_00"
$ws.Cells.Item(201, 2).Value = 34
$ws.Cells.Item(201, 3).Value = 0
$ws.Cells.Item(201, 4).Value = 0
$ws.Cells.Item(201, 5).Value = 0
$ws.Cells.Item(201, 6).Value = 39
$ws.Cells.Item(201, 7).Value = 0
$ws.Cells.Item(201, 8).Value = 0
$ws.Cells.Item(201, 9).Value = 0

$ws.Cells.Item(202, 1).Value = "20_06"
$ws.Cells.Item(202, 2).Value = 0
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 0
$ws.Cells.Item(202, 5).Value = 0
$ws.Cells.Item(202, 6).Value = 35
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0
$ws.Cells.Item(202, 9).Value = 0
